# 2021TPI_BG_JdB.docx -- "avancée de la documentation + correction nom fichier"
#
# 1) The "J9" heading paragraph currently carries its text split across six
#    runs ("J9" / " " / "Mardi" / " 1" / "8" / " mai 2021"); collapse it to a
#    single run with the same visible text "J9 Mardi 18 mai 2021".
# 2) Append three new log entries after the existing "07:30 : ..." paragraph
#    in the J9 "Déroulement" section, and move the (normally hidden)
#    "_GoBack" bookmark so it once again sits right after the very last
#    character typed ("11:33 : Je sauvegarde"), matching real Word's
#    behaviour of tracking the last edit position.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1 -- merge the "J9 ... mai 2021" heading runs into a single run
# ---------------------------------------------------------------------

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("J9") -and $cand.Style.NameLocal -eq "Heading 2") {
        $headingIndex = $i
        $found = $true
        break
    }
}

if ($found) {
    $heading = $d.Paragraphs.Item($headingIndex)
    # Range excluding the trailing paragraph mark.
    $headRange = $d.Range($heading.Range.Start, $heading.Range.End - 1)
    $plainText = $headRange.Text

    # Clear then retype so the run-splitting collapses to a single run even
    # though the visible characters end up identical to before.
    $headRange.Text = ""
    $heading2 = $d.Paragraphs.Item($headingIndex)
    $headRange2 = $d.Range($heading2.Range.Start, $heading2.Range.End - 1)
    $headRange2.Text = $plainText
}

# ---------------------------------------------------------------------
# Part 2 -- add the new "Déroulement" entries and relocate "_GoBack"
# ---------------------------------------------------------------------

$bm = $d.Bookmarks("_GoBack")
$startPos = $bm.Start
$insertPoint = $d.Range($startPos, $startPos)

$entry1 = "08:00 : Je discute avec M.Garcia de la méthode recommandée pour sortir le code source. La méthode proposée est avec Notepad++."
$entry2 = "08:10 : Je reprends la documentation à la navigation du site"
$entry3 = "11:33 : Je sauvegarde"

# Insert the three entries separated by a placeholder marker, right at the
# bookmark's current position -- this keeps growing the SAME paragraph as
# the bookmark rides along with it. A rarely-used Private-Use-Area
# character is used as the placeholder (instead of a literal CR) so a
# later scoped Find/Replace can turn just these marks into real paragraph
# breaks without touching unrelated manual line breaks (`<w:br/>`, which
# Word's Range.Text also represents as a vertical-tab) that already exist
# elsewhere in the document.
$ph = [char]0xE000
$combined = $ph + $entry1 + $ph + $entry2 + $ph + $entry3
$insertPoint.InsertAfter($combined)
$endPos = $startPos + $combined.Length

# Turn the placeholder marks into paragraph breaks, scoped to only the
# text we just inserted. Once this runs, the bookmark naturally ends up
# glued to the end of "Je sauvegarde", exactly like Word leaves "_GoBack"
# at the last character edited.
$scoped = $d.Range($startPos, $endPos)
$scoped.Find.ClearFormatting()
$scoped.Find.Replacement.ClearFormatting()
$scoped.Find.Execute($ph, $false, $false, $false, $false, $false, $true, 1, $false, "`r", 2)
